$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40 (ALC)
$ws.Range("H40").Value = 6332.3335
$ws.Range("I40").Value = 4498.5
$ws.Range("J40").Value = 10000
$ws.Range("K40").Value = 4498.5
$ws.Range("L40").Value = 10000
$ws.Range("M40").Value = -4323.5
$ws.Range("N40").Value = -10350

# Row 51 (ALC)
$ws.Range("H51").Value = 14624.5
$ws.Range("I51").Value = 14624.5
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 14624.5
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -14140.5

# Row 62 (ALC)
$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -4376

# Row 65 (ALC)
$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -21880

# Row 100 (ALC)
$ws.Range("H100").Value = 1569.6
$ws.Range("I100").Value = 1466.6666
$ws.Range("J100").Value = 1724
$ws.Range("K100").Value = 1466.6666
$ws.Range("L100").Value = 1724
$ws.Range("M100").Value = -925.6666
$ws.Range("N100").Value = -2806

$ws = $wb.Worksheets.Item("ARM")
# Row 31 (ARM)
$ws.Range("H31").Value = 14999
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 14999
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 14999
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -15587

# Row 32 (ARM)
$ws.Range("H32").Value = 3933.2563
$ws.Range("I32").Value = 3510.4473
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 3510.4473
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -3223.4473

# Row 45 (ARM)
$ws.Range("H45").Value = 3499.75
$ws.Range("I45").Value = 1999
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 1999
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1622

# Row 61 (ARM)
$ws.Range("H61").Value = 2406.7222
$ws.Range("I61").Value = 1146.5454
$ws.Range("J61").Value = 4387
$ws.Range("K61").Value = 1146.5454
$ws.Range("L61").Value = 4387
$ws.Range("M61").Value = -934.5454
$ws.Range("N61").Value = -4811

# Row 102 (ARM)
$ws.Range("H102").Value = 1957
$ws.Range("I102").Value = 1957
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1957
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -335
$ws.Range("N102").ClearContents()

# Row 136 (ARM)
$ws.Range("H136").Value = 2406.7222
$ws.Range("I136").Value = 1146.5454
$ws.Range("J136").Value = 4387
$ws.Range("K136").Value = 3439.6362
$ws.Range("L136").Value = 13161
$ws.Range("M136").Value = -889.6361999999999
$ws.Range("N136").Value = -18261

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (BSM)
$ws.Range("H20").Value = 3319.9333
$ws.Range("I20").Value = 3190
$ws.Range("J20").Value = 3677.25
$ws.Range("K20").Value = 3190
$ws.Range("L20").Value = 3677.25
$ws.Range("M20").Value = -2943
$ws.Range("N20").Value = -4171.25

# Row 99 (BSM)
$ws.Range("H99").Value = 737.5
$ws.Range("I99").Value = 737.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 737.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 760.5

# Row 105 (BSM)
$ws.Range("H105").Value = 3596.6667
$ws.Range("I105").Value = 2515.8
$ws.Range("J105").Value = 9001
$ws.Range("K105").Value = 2515.8
$ws.Range("L105").Value = 9001
$ws.Range("M105").Value = -768.8000000000002
$ws.Range("N105").Value = -12495

# Row 134 (BSM)
$ws.Range("H134").Value = 9612.6875
$ws.Range("I134").Value = 3350.5
$ws.Range("J134").Value = 15874.875
$ws.Range("K134").Value = 10051.5
$ws.Range("L134").Value = 47624.625
$ws.Range("M134").Value = -7516.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3323.5789
$ws.Range("I31").Value = 1250
$ws.Range("J31").Value = 4533.1665
$ws.Range("K31").Value = 1250
$ws.Range("L31").Value = 4533.1665
$ws.Range("M31").Value = -955
$ws.Range("N31").Value = -5123.1665

# Row 34 (CRP)
$ws.Range("H34").Value = 3323.5789
$ws.Range("I34").Value = 1250
$ws.Range("J34").Value = 4533.1665
$ws.Range("K34").Value = 1250
$ws.Range("L34").Value = 4533.1665
$ws.Range("M34").Value = -1048
$ws.Range("N34").Value = -4937.1665

# Row 105 (CRP)
$ws.Range("H105").Value = 3086.6667
$ws.Range("I105").Value = 3086.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3086.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -1339.6667

# Row 112 (CRP)
$ws.Range("H112").Value = 15000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 15000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 15000
$ws.Range("N112").Value = -17954

# Row 132 (CRP)
$ws.Range("H132").Value = 4798.8
$ws.Range("I132").Value = 4498.5
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 13495.5
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -10965.5

# Row 134 (CRP)
$ws.Range("H134").Value = 4999.4
$ws.Range("I134").Value = 4999.4
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14998.2
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12463.2

$ws = $wb.Worksheets.Item("CUL")
# Row 4 (CUL)
$ws.Range("H4").Value = 4494444.5
$ws.Range("I4").Value = 83614
$ws.Range("J4").Value = 22137766
$ws.Range("K4").Value = 250842
$ws.Range("L4").Value = 66413298
$ws.Range("M4").Value = -250730
$ws.Range("N4").Value = -66413522

# Row 5 (CUL)
$ws.Range("H5").Value = 1237.875
$ws.Range("I5").Value = 1620.8
$ws.Range("J5").Value = 599.6667
$ws.Range("K5").Value = 4862.4
$ws.Range("L5").Value = 1799.0001
$ws.Range("M5").Value = -4750.4
$ws.Range("N5").Value = -2023.0001

# Row 33 (CUL)
$ws.Range("H33").Value = 543.5
$ws.Range("I33").Value = 698
$ws.Range("J33").Value = 389
$ws.Range("K33").Value = 4188
$ws.Range("L33").Value = 2334
$ws.Range("M33").Value = -3905

# Row 60 (CUL)
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()

# Row 70 (CUL)
$ws.Range("H70").Value = 1770.6666
$ws.Range("I70").Value = 256
$ws.Range("J70").Value = 4800
$ws.Range("K70").Value = 768
$ws.Range("L70").Value = 14400
$ws.Range("M70").Value = -453
$ws.Range("N70").Value = -15030

# Row 73 (CUL)
$ws.Range("H73").Value = 1770.6666
$ws.Range("I73").Value = 256
$ws.Range("J73").Value = 4800
$ws.Range("K73").Value = 768
$ws.Range("L73").Value = 14400
$ws.Range("M73").Value = 324
$ws.Range("N73").Value = -16584

# Row 131 (CUL)
$ws.Range("H131").Value = 2103
$ws.Range("I131").Value = 20
$ws.Range("J131").Value = 2276.5833
$ws.Range("K131").Value = 60
$ws.Range("L131").Value = 6829.749899999999
$ws.Range("M131").Value = 4980
$ws.Range("N131").Value = -16909.7499

# Row 135 (CUL)
$ws.Range("H135").Value = 1237.875
$ws.Range("I135").Value = 1620.8
$ws.Range("J135").Value = 599.6667
$ws.Range("K135").Value = 14587.2
$ws.Range("L135").Value = 5397.0003
$ws.Range("M135").Value = -12052.2
$ws.Range("N135").Value = -10467.0003

$ws = $wb.Worksheets.Item("GSM")
# Row 97 (GSM)
$ws.Range("H97").Value = 1302
$ws.Range("I97").Value = 992.8
$ws.Range("J97").Value = 2075
$ws.Range("K97").Value = 992.8
$ws.Range("L97").Value = 2075
$ws.Range("M97").Value = -496.8

# Row 113 (GSM)
$ws.Range("H113").Value = 2910.75
$ws.Range("I113").Value = 2910.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2910.75
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -740.75

# Row 126 (GSM)
$ws.Range("H126").Value = 3671
$ws.Range("I126").Value = 1599
$ws.Range("J126").Value = 4707
$ws.Range("K126").Value = 4797
$ws.Range("L126").Value = 14121
$ws.Range("M126").Value = -2327

# Row 132 (GSM)
$ws.Range("H132").Value = 3223.9
$ws.Range("I132").Value = 2780.125
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8340.375
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5810.375

$ws = $wb.Worksheets.Item("LTW")
# Row 22 (LTW)
$ws.Range("H22").Value = 7180
$ws.Range("I22").Value = 1266.6666
$ws.Range("J22").Value = 9714.286
$ws.Range("K22").Value = 1266.6666
$ws.Range("L22").Value = 9714.286
$ws.Range("M22").Value = -971.6666
$ws.Range("N22").Value = -10304.286

# Row 27 (LTW)
$ws.Range("H27").Value = 7180
$ws.Range("I27").Value = 1266.6666
$ws.Range("J27").Value = 9714.286
$ws.Range("K27").Value = 1266.6666
$ws.Range("L27").Value = 9714.286
$ws.Range("M27").Value = -1159.6666
$ws.Range("N27").Value = -9928.286

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (WVR)
$ws.Range("I62").Value = 6000
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -5376
$ws.Range("N62").Value = -7248

# Row 65 (WVR)
$ws.Range("I65").Value = 6000
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -26880
$ws.Range("N65").Value = -36240

# Row 136 (WVR)
$ws.Range("H136").Value = 2757.7144
$ws.Range("I136").Value = 2560.8
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 7682.400000000001
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -5132.400000000001
